$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("#i.-sone",    "I. Sone"),
    @("#i.sone",     "I.Sone"),
    @("#iii.sone",   "III.Sone"),
    @("#ii.-sone",   "II. Sone"),
    @("#ii-sone",    "II Sone"),
    @("#iii.-sone",  "III. Sone"),
    @("#i.-zone",    "I. Zone"),
    @("#moeder",     "Moeder"),
    @("#vader",      "Vader"),
    @("#i.zone",     "I.Zone"),
    @("#i-zone",     "I Zone"),
    @("#ii.sone",    "II.Sone"),
    @("#ii.zone",    "II.Zone"),
    @("#[...]-sone", "[...] Sone"),
    @("#roel",       "Roel"),
    @("#sone",       "Sone"),
    @("#zone",       "Zone"),
    @("#may",        "May")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = ""
}
